$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# Footer Pearson logo pictures: rename from "image1.png" to "image2.png"
for ($i = 1; $i -le 2; $i++) {
    $ftr = $sec.Footers.Item($i)
    $shp = $ftr.Range.InlineShapes.Item(1)
    $shp.Name = "image2.png"
}

# Header BTEC logo pictures: rename from "image2.jpg" to "image1.jpg"
for ($i = 1; $i -le 2; $i++) {
    $hdr = $sec.Headers.Item($i)
    $shp = $hdr.Range.InlineShapes.Item(1)
    $shp.Name = "image1.jpg"
}

Write-Output "Renamed Pearson footer logos to image2.png and BTEC header logos to image1.jpg"
